$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 81
$ws1.Range("F6").Value = 729
$ws1.Range("F8").Value = 229
$ws1.Range("F17").Value = 799
$ws1.Range("F27").Value = 950
$ws1.Range("F29").Value = 187

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 26
$ws2.Range("F4").Value = 990
$ws2.Range("F5").Value = 990
$ws2.Range("F17").Value = 961
$ws2.Range("F31").Value = 77

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2369
$ws3.Range("F6").Value = 962
$ws3.Range("F9").Value = 1204
$ws3.Range("F10").Value = 311

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2369
$ws4.Range("F7").Value = 26
$ws4.Range("F8").Value = 962
$ws4.Range("F9").Value = 1204
$ws4.Range("F10").Value = 311
$ws4.Range("F13").Value = 81
$ws4.Range("F14").Value = 729
$ws4.Range("F17").Value = 229
$ws4.Range("F21").Value = 990
$ws4.Range("F26").Value = 799
$ws4.Range("F34").Value = 950
$ws4.Range("F38").Value = 187
